# Auto-applies the per-cell numeric corrections described in the commit diff.
# Columns H:N hold market-derived pricing/profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) that were refreshed
# by the scheduled pricing-data runner; row/label columns A:G are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 49
  $ws.Range("H49").Value = 483.375
  $ws.Range("I49").Value = 166.75
  $ws.Range("J49").Value = 800
  $ws.Range("K49").Value = 500.25
  $ws.Range("L49").Value = 2400
  $ws.Range("M49").Value = -364.25
  $ws.Range("N49").Value = -2672
  # Row 69
  $ws.Range("H69").Value = 3400
  $ws.Range("I69").Value = 3400
  $ws.Range("J69").Value = 3400
  $ws.Range("K69").Value = 10200
  $ws.Range("L69").Value = 10200
  $ws.Range("M69").Value = -9326
  $ws.Range("N69").Value = -11948
  # Row 72
  $ws.Range("H72").Value = 3400
  $ws.Range("I72").Value = 3400
  $ws.Range("J72").Value = 3400
  $ws.Range("K72").Value = 30600
  $ws.Range("L72").Value = 30600
  $ws.Range("M72").Value = -26232
  $ws.Range("N72").Value = -39336
  # Row 129
  $ws.Range("H129").Value = 2410.8657
  $ws.Range("I129").Value = 690.8
  $ws.Range("J129").Value = 2549.5806
  $ws.Range("K129").Value = 2072.4
  $ws.Range("L129").Value = 7648.7418
  $ws.Range("M129").Value = 2927.6
  $ws.Range("N129").Value = -17648.7418
  # Row 138
  $ws.Range("H138").Value = 2508.675
  $ws.Range("I138").Value = 2640.3684
  $ws.Range("J138").Value = 2389.524
  $ws.Range("K138").Value = 7921.1052
  $ws.Range("L138").Value = 7168.572
  $ws.Range("M138").Value = -2781.1052
  $ws.Range("N138").Value = -17448.572

$ws = $wb.Worksheets.Item("ARM")
  # Row 32
  $ws.Range("H32").Value = 15976.103
  $ws.Range("I32").Value = 16836.625
  $ws.Range("J32").Value = 5649.8335
  $ws.Range("K32").Value = 16836.625
  $ws.Range("L32").Value = 5649.8335
  $ws.Range("M32").Value = -16549.625
  $ws.Range("N32").Value = -6223.8335
  # Row 132
  $ws.Range("H132").Value = 4165.1963
  $ws.Range("I132").Value = 4797.6787
  $ws.Range("J132").Value = 3395.2173
  $ws.Range("K132").Value = 14393.0361
  $ws.Range("L132").Value = 10185.6519
  $ws.Range("M132").Value = -11863.0361
  $ws.Range("N132").Value = -15245.6519

$ws = $wb.Worksheets.Item("BSM")
  # Row 86
  $ws.Range("H86").Value = 1500
  $ws.Range("I86").Value = 1500
  $ws.Range("J86").Value = 0
  $ws.Range("K86").Value = 1500
  $ws.Range("L86").Value = 0
  $ws.Range("M86").Value = -377
  $ws.Range("N86").Value = $null
  # Row 89
  $ws.Range("H89").Value = 1500
  $ws.Range("I89").Value = 1500
  $ws.Range("J89").Value = 0
  $ws.Range("K89").Value = 7500
  $ws.Range("L89").Value = 0
  $ws.Range("M89").Value = -1884
  $ws.Range("N89").Value = $null
  # Row 94
  $ws.Range("H94").Value = 450.2
  $ws.Range("I94").Value = 448.125
  $ws.Range("J94").Value = 500
  $ws.Range("K94").Value = 448.125
  $ws.Range("L94").Value = 500
  $ws.Range("M94").Value = 2.875
  $ws.Range("N94").Value = -1402
  # Row 107
  $ws.Range("H107").Value = 1610.1578
  $ws.Range("I107").Value = 1051
  $ws.Range("J107").Value = 2379
  $ws.Range("K107").Value = 1051
  $ws.Range("L107").Value = 2379
  $ws.Range("M107").Value = 869
  $ws.Range("N107").Value = -6219

$ws = $wb.Worksheets.Item("CRP")
  # Row 16
  $ws.Range("H16").Value = 2280
  $ws.Range("I16").Value = 1133.3334
  $ws.Range("J16").Value = 4000
  $ws.Range("K16").Value = 1133.3334
  $ws.Range("L16").Value = 4000
  $ws.Range("M16").Value = -846.3334
  $ws.Range("N16").Value = -4574
  # Row 31
  $ws.Range("H31").Value = 3573286.8
  $ws.Range("I31").Value = 2470.318
  $ws.Range("J31").Value = 5883815
  $ws.Range("K31").Value = 2470.318
  $ws.Range("L31").Value = 5883815
  $ws.Range("M31").Value = -2175.318
  $ws.Range("N31").Value = -5884405
  # Row 34
  $ws.Range("H34").Value = 3573286.8
  $ws.Range("I34").Value = 2470.318
  $ws.Range("J34").Value = 5883815
  $ws.Range("K34").Value = 2470.318
  $ws.Range("L34").Value = 5883815
  $ws.Range("M34").Value = -2268.318
  $ws.Range("N34").Value = -5884219
  # Row 99
  $ws.Range("H99").Value = 2484.8948
  $ws.Range("I99").Value = 2033.2222
  $ws.Range("K99").Value = 2033.2222
  $ws.Range("M99").Value = -535.2221999999999
  # Row 113
  $ws.Range("H113").Value = 2280
  $ws.Range("I113").Value = 1133.3334
  $ws.Range("J113").Value = 4000
  $ws.Range("K113").Value = 1133.3334
  $ws.Range("L113").Value = 4000
  $ws.Range("M113").Value = 1036.6666
  $ws.Range("N113").Value = -8340
  # Row 126
  $ws.Range("H126").Value = 2484.8948
  $ws.Range("I126").Value = 2033.2222
  $ws.Range("K126").Value = 6099.6666
  $ws.Range("M126").Value = -3629.6666
  # Row 141
  $ws.Range("H141").Value = 52079.168
  $ws.Range("J141").Value = 52079.168
  $ws.Range("L141").Value = 52079.168
  $ws.Range("N141").Value = -62439.168

$ws = $wb.Worksheets.Item("CUL")
  # Row 104
  $ws.Range("H104").Value = 251
  $ws.Range("I104").Value = 251
  $ws.Range("J104").Value = 0
  $ws.Range("K104").Value = 753
  $ws.Range("L104").Value = 0
  $ws.Range("M104").Value = 1868
  $ws.Range("N104").Value = $null
  # Row 131
  $ws.Range("H131").Value = 2317807
  $ws.Range("J131").Value = 4116254
  $ws.Range("L131").Value = 12348762
  $ws.Range("N131").Value = -12358842

$ws = $wb.Worksheets.Item("GSM")
  # Row 107
  $ws.Range("H107").Value = 542.7241
  $ws.Range("I107").Value = 510.5909
  $ws.Range("J107").Value = 643.7143
  $ws.Range("K107").Value = 510.5909
  $ws.Range("L107").Value = 643.7143
  $ws.Range("M107").Value = 1409.4091
  $ws.Range("N107").Value = -4483.7143
  # Row 113
  $ws.Range("H113").Value = 13159740
  $ws.Range("I113").Value = 35715316
  $ws.Range("J113").Value = 2320.8333
  $ws.Range("K113").Value = 35715316
  $ws.Range("L113").Value = 2320.8333
  $ws.Range("M113").Value = -35713146
  $ws.Range("N113").Value = -6660.8333
  # Row 132
  $ws.Range("H132").Value = 63659.09
  $ws.Range("I132").Value = 102461.9
  $ws.Range("J132").Value = 3962.4614
  $ws.Range("K132").Value = 307385.7
  $ws.Range("L132").Value = 11887.3842
  $ws.Range("M132").Value = -304855.7
  $ws.Range("N132").Value = -16947.3842

$ws = $wb.Worksheets.Item("LTW")
  # Row 22
  $ws.Range("H22").Value = 601.55554
  $ws.Range("I22").Value = 550
  $ws.Range("J22").Value = 642.8
  $ws.Range("K22").Value = 550
  $ws.Range("L22").Value = 642.8
  $ws.Range("M22").Value = -255
  $ws.Range("N22").Value = -1232.8
  # Row 27
  $ws.Range("H27").Value = 601.55554
  $ws.Range("I27").Value = 550
  $ws.Range("J27").Value = 642.8
  $ws.Range("K27").Value = 550
  $ws.Range("L27").Value = 642.8
  $ws.Range("M27").Value = -443
  $ws.Range("N27").Value = -856.8
  # Row 46
  $ws.Range("H46").Value = 1831.5
  $ws.Range("I46").Value = 1495.7858
  $ws.Range("J46").Value = 2125.25
  $ws.Range("K46").Value = 1495.7858
  $ws.Range("L46").Value = 2125.25
  $ws.Range("M46").Value = -1307.7858
  $ws.Range("N46").Value = -2501.25
  # Row 100
  $ws.Range("H100").Value = 3399.4
  $ws.Range("I100").Value = 3000.75
  $ws.Range("K100").Value = 3000.75
  $ws.Range("M100").Value = -2459.75

$ws = $wb.Worksheets.Item("WVR")
  # Row 81
  $ws.Range("H81").Value = 2600
  $ws.Range("I81").Value = 2500
  $ws.Range("J81").Value = 2700
  $ws.Range("K81").Value = 5000
  $ws.Range("L81").Value = 5400
  $ws.Range("M81").Value = -3939
  $ws.Range("N81").Value = -7522
  # Row 84
  $ws.Range("H84").Value = 2600
  $ws.Range("I84").Value = 2500
  $ws.Range("J84").Value = 2700
  $ws.Range("K84").Value = 25000
  $ws.Range("L84").Value = 27000
  $ws.Range("M84").Value = -19696
  $ws.Range("N84").Value = -37608
  # Row 122
  $ws.Range("H122").Value = 3168
  $ws.Range("I122").Value = 3210.889
  $ws.Range("J122").Value = 2975
  $ws.Range("K122").Value = 9632.667000000001
  $ws.Range("L122").Value = 8925
  $ws.Range("M122").Value = -7182.667000000001
  $ws.Range("N122").Value = -13825
  # Row 136
  $ws.Range("H136").Value = 1749.2413
  $ws.Range("I136").Value = 1774.125
  $ws.Range("J136").Value = 1629.8
  $ws.Range("K136").Value = 5322.375
  $ws.Range("L136").Value = 4889.4
  $ws.Range("M136").Value = -2772.375
  $ws.Range("N136").Value = -9989.4
